# daily auto push: 2025-09-30 13:37 UTC
# Append the latest reading (2025/09/30, 20:00, rank 11) as a new row at the
# bottom of the log table on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 41

# Column A stores the date as literal text (e.g. "2025/09/22") just like
# every other row above it, not as a real Excel date serial number. Force
# text interpretation via NumberFormat, then reset the style back to the
# workbook default so the new row doesn't pick up an extra explicit style.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025/09/30"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = "火"
$ws.Cells.Item($row, 3).Value = 20
$ws.Cells.Item($row, 4).Value = 11
